# Update "想去人数" (interested-people count) figures in column F
# for both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 471
$ws1.Range("F7").Value = 221
$ws1.Range("F8").Value = 14589
$ws1.Range("F9").Value = 167
$ws1.Range("F10").Value = 126
$ws1.Range("F11").Value = 5830
$ws1.Range("F12").Value = 594
$ws1.Range("F13").Value = 77
$ws1.Range("F15").Value = 66
$ws1.Range("F18").Value = 84
$ws1.Range("F22").Value = 59
$ws1.Range("F23").Value = 10611
$ws1.Range("F26").Value = 86

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 471
$ws4.Range("F8").Value = 221
$ws4.Range("F9").Value = 14589
$ws4.Range("F10").Value = 167
$ws4.Range("F11").Value = 126
$ws4.Range("F12").Value = 5830
$ws4.Range("F13").Value = 594
$ws4.Range("F14").Value = 77
$ws4.Range("F16").Value = 66
$ws4.Range("F19").Value = 84
$ws4.Range("F23").Value = 59
$ws4.Range("F25").Value = 10611
$ws4.Range("F28").Value = 86
